$d = $word.ActiveDocument

# --- Paragraph 1 ("**ID__AFFARS_pgi_5301_topic_26__ID**" placeholder line) ---
$p1 = $d.Paragraphs(1)

# Add a paragraph border (5-twip space on every side) and widen the left
# indent from 120 -> 225 twips (6pt -> 11.25pt).
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# The paragraph currently holds two runs: the ID placeholder text, then a
# trailing " " run. Work out their boundaries from the paragraph range
# (Range.End sits just past the trailing paragraph mark).
$oldId  = "**ID__AFFARS_pgi_5301_topic_26__ID**"
$newId  = "**ID__AFFARS_AFMC_PGI_5301_601B__ID**"
$pStart = $p1.Range.Start
$idEnd  = $pStart + $oldId.Length
$pMarkStart = $p1.Range.End - 1

# Remove the trailing space-only run first so the ID run isn't re-merged
# with it, then update the ID run's text in place.
$spaceRun = $d.Range($idEnd, $pMarkStart)
$spaceRun.Delete()

$idRun = $d.Range($pStart, $idEnd)
$idRun.Text = $newId
